$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "27.393.46"
Set-TextValue $ws.Range("E2") "  +3.72%  "
Set-TextValue $ws.Range("D3") "1.841.65"
Set-TextValue $ws.Range("E3") "  +4.02%  "
Set-TextValue $ws.Range("D4") "1.024"
Set-TextValue $ws.Range("E4") "  +2.94%  "
Set-TextValue $ws.Range("D5") "318.67"
Set-TextValue $ws.Range("E5") "  +4.37%  "
Set-TextValue $ws.Range("D6") "1.023"
Set-TextValue $ws.Range("E6") "  +2.69%  "
Set-TextValue $ws.Range("D7") "0.4350"
Set-TextValue $ws.Range("D8") "0.3722"
Set-TextValue $ws.Range("E8") "  +2.23%  "
Set-TextValue $ws.Range("D9") "0.07328"
Set-TextValue $ws.Range("E9") "  +2.15%  "
Set-TextValue $ws.Range("D10") "0.8760"
Set-TextValue $ws.Range("E10") "  +3.25%  "
Set-TextValue $ws.Range("D11") "21.35"
Set-TextValue $ws.Range("E11") "  +4.78%  "
Set-TextValue $ws.Range("D12") "2.016.73"
Set-TextValue $ws.Range("E12") "  +15.14%  "
Set-TextValue $ws.Range("D13") "5.479"
Set-TextValue $ws.Range("E13") "  +4.54%  "
Set-TextValue $ws.Range("D14") "6.678"
Set-TextValue $ws.Range("E14") "  +3.86%  "
Set-TextValue $ws.Range("D15") "0.07159"
Set-TextValue $ws.Range("E15") "  +4.63%  "
Set-TextValue $ws.Range("D16") "82.11"
Set-TextValue $ws.Range("E16") "  +4.12%  "
Set-TextValue $ws.Range("E17") "  +2.88%  "
Set-TextValue $ws.Range("E18") "  +3.74%  "
Set-TextValue $ws.Range("D19") "1.020"
Set-TextValue $ws.Range("E19") "  +2.46%  "
Set-TextValue $ws.Range("D20") "15.43"
Set-TextValue $ws.Range("E20") "  +2.84%  "
Set-TextValue $ws.Range("D21") "27.388.20"
Set-TextValue $ws.Range("E21") "  +3.77%  "
Set-TextValue $ws.Range("D22") "5.241"
Set-TextValue $ws.Range("E22") "  +2.85%  "
Set-TextValue $ws.Range("D23") "11.11"
Set-TextValue $ws.Range("E23") "  +0.09%  "
Set-TextValue $ws.Range("D24") "2.222.30"
Set-TextValue $ws.Range("E24") "  +12.31%  "
Set-TextValue $ws.Range("D25") "156.71"
Set-TextValue $ws.Range("E25") "  +3.21%  "
Set-TextValue $ws.Range("D26") "1.907"
Set-TextValue $ws.Range("E26") "  +2.32%  "
Set-TextValue $ws.Range("D27") "18.53"
Set-TextValue $ws.Range("E27") "  +2.76%  "
Set-TextValue $ws.Range("D28") "5.274"
Set-TextValue $ws.Range("D29") "1.928"
Set-TextValue $ws.Range("E29") "  +6.16%  "
Set-TextValue $ws.Range("D30") "115.50"
Set-TextValue $ws.Range("E30") "  +1.57%  "
Set-TextValue $ws.Range("D31") "0.09008"
Set-TextValue $ws.Range("E31") "  +0.76%  "
Set-TextValue $ws.Range("D32") "1.200"
Set-TextValue $ws.Range("E32") "  +6.79%  "
Set-TextValue $ws.Range("D33") "0.7602"
Set-TextValue $ws.Range("E33") "  +4.03%  "
Set-TextValue $ws.Range("D34") "4.457"
Set-TextValue $ws.Range("E34") "  +3.06%  "
Set-TextValue $ws.Range("D35") "2.845"
Set-TextValue $ws.Range("E35") "  +4.29%  "
Set-TextValue $ws.Range("D36") "1.025"
Set-TextValue $ws.Range("E36") "  +2.91%  "
Set-TextValue $ws.Range("E37") "  +4.85%  "
Set-TextValue $ws.Range("D38") "0.01953"
Set-TextValue $ws.Range("E38") "  +3.84%  "
Set-TextValue $ws.Range("D39") "0.05265"
Set-TextValue $ws.Range("E39") "  +2.20%  "
Set-TextValue $ws.Range("D42") "0.1662"
Set-TextValue $ws.Range("E42") "  +3.11%  "
Set-TextValue $ws.Range("D43") "6.528"
Set-TextValue $ws.Range("E43") "  +3.38%  "
Set-TextValue $ws.Range("D44") "8.468"
Set-TextValue $ws.Range("E44") "  +5.40%  "
Set-TextValue $ws.Range("D45") "108.11"
Set-TextValue $ws.Range("E45") "  +3.02%  "
Set-TextValue $ws.Range("D46") "10.55"
Set-TextValue $ws.Range("E46") "  +3.88%  "
Set-TextValue $ws.Range("D47") "1.026"
Set-TextValue $ws.Range("E47") "  +3.04%  "
Set-TextValue $ws.Range("D48") "0.4632"
Set-TextValue $ws.Range("E48") "  +3.12%  "
Set-TextValue $ws.Range("D49") "1.667"
Set-TextValue $ws.Range("E49") "  +2.99%  "
Set-TextValue $ws.Range("D50") "1.896"
Set-TextValue $ws.Range("E50") "  +8.97%  "
Set-TextValue $ws.Range("D51") "0.06287"
Set-TextValue $ws.Range("E51") "  +1.52%  "

# Row 40/41 swap (TheSandbox <-> MXToken)
Set-TextValue $ws.Range("B40") "MXToken"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D40") "2.809"
Set-TextValue $ws.Range("E40") "  +8.84%  "
Set-TextValue $ws.Range("B41") "TheSandbox"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D41") "0.5154"
Set-TextValue $ws.Range("E41") "  +4.53%  "
